# rerun of models including new transit access feature
# Updates regression coefficients/p-values on both worksheets and renames the
# "DistCenter" parameter row to "DistCenter_pc" on the All_model_short sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("All_model_short")
$ws2 = $wb.Worksheets.Item("All_model_short_table")

# --- Sheet "All_model_short" (columns A:D) ---
$ws1.Range("B2").Value = -0.46023953223507202
$ws1.Range("C2").Value = [double]"1.6980652489026601E-20"
$ws1.Range("B3").Value = [double]"5.5707971072500005E-4"
$ws1.Range("B4").Value = 0.26762141147343699
$ws1.Range("C4").Value = [double]"3.06692592906174E-72"
$ws1.Range("B5").Value = [double]"1.4438622467591301E-2"
$ws1.Range("C5").Value = [double]"1.5167297187798899E-34"
$ws1.Range("B6").Value = [double]"1.7964720611974301E-2"
$ws1.Range("C6").Value = 0.584482422169879

# Row 11: parameter renamed DistCenter -> DistCenter_pc, plus new coefficient/p
$ws1.Range("A11").Value = "DistCenter_pc"
$ws1.Range("B11").Value = [double]"2.6592338875686002E-3"
$ws1.Range("C11").Value = [double]"4.2523353552762901E-10"
$ws1.Range("C11").NumberFormat = "0.00E+00"

$ws1.Range("B12").Value = [double]"-1.26940107793429E-2"
$ws1.Range("C12").Value = 0.22159900944383601
$ws1.Range("B13").Value = [double]"-7.2404768539832999E-3"
$ws1.Range("C13").Value = [double]"4.1602598350801603E-119"
$ws1.Range("B15").Value = [double]"-1.1497969342928E-3"
$ws1.Range("C15").Value = 0.39047027816843799
$ws1.Range("B16").Value = [double]"-3.5020744151725998E-3"
$ws1.Range("C16").Value = [double]"1.7464723909840001E-4"
$ws1.Range("B17").Value = 0.25676778313603299
$ws1.Range("C17").Value = 0.119577024469868
$ws1.Range("B18").Value = [double]"-9.9248980255405E-3"
$ws1.Range("C18").Value = [double]"3.1905094989472801E-28"
$ws1.Range("B19").Value = 0.76012747802402003
$ws1.Range("C19").Value = [double]"2.2370683587341199E-11"
$ws1.Range("B20").Value = -0.63594752240793095
$ws1.Range("C20").Value = [double]"1.8968741510654099E-5"

# Match the selection left behind in the saved file
$ws1.Range("B15:C20").Select() | Out-Null

# --- Sheet "All_model_short_table" (columns B:D) ---
$ws2.Range("C2").Value = -0.46023953223507202
$ws2.Range("D2").Value = [double]"1.6980652489026601E-20"
$ws2.Range("C3").Value = [double]"5.5707971072500005E-4"
$ws2.Range("C4").Value = 0.26762141147343699
$ws2.Range("D4").Value = [double]"3.06692592906174E-72"
$ws2.Range("C5").Value = [double]"1.4438622467591301E-2"
$ws2.Range("D5").Value = [double]"1.5167297187798899E-34"
$ws2.Range("C6").Value = [double]"1.7964720611974301E-2"
$ws2.Range("D6").Value = 0.584482422169879
$ws2.Range("C7").Value = [double]"2.6592338875686002E-3"
$ws2.Range("D7").Value = [double]"4.2523353552762901E-10"
$ws2.Range("C8").Value = [double]"-1.26940107793429E-2"
$ws2.Range("D8").Value = 0.22159900944383601
$ws2.Range("C9").Value = [double]"-7.2404768539832999E-3"
$ws2.Range("D9").Value = [double]"4.1602598350801603E-119"
$ws2.Range("C10").Value = [double]"-1.1497969342928E-3"
$ws2.Range("D10").Value = 0.39047027816843799
$ws2.Range("C11").Value = [double]"-3.5020744151725998E-3"
$ws2.Range("D11").Value = [double]"1.7464723909840001E-4"
$ws2.Range("C12").Value = 0.25676778313603299
$ws2.Range("D12").Value = 0.119577024469868
$ws2.Range("C13").Value = [double]"-9.9248980255405E-3"
$ws2.Range("D13").Value = [double]"3.1905094989472801E-28"
$ws2.Range("C14").Value = 0.76012747802402003
$ws2.Range("D14").Value = [double]"2.2370683587341199E-11"
$ws2.Range("C15").Value = -0.63594752240793095
$ws2.Range("D15").Value = [double]"1.8968741510654099E-5"

# Keep "All_model_short_table" the active sheet and match its saved selection
$ws2.Activate()
$ws2.Range("C10:D15").Select() | Out-Null
